$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, pushing existing rows 68-97 down to 69-98.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new data record.
$ws.Range("A68").Value = 11
$ws.Range("B68").Value = "Vega Monumental Concepción"
$ws.Range("C68").Value = "Bíobío"
$ws.Range("D68").Value = 44572
$ws.Range("E68").Value = 8
$ws.Range("F68").Value = "Fruta"
$ws.Range("G68").Value = 100109
$ws.Range("H68").Value = "Uva"
$ws.Range("I68").Value = 100109001
$ws.Range("J68").Value = "Uva"
$ws.Range("K68").Value = "Superior Seedless"
$ws.Range("L68").Value = "Primera"
$ws.Range("M68").Value = 100
$ws.Range("N68").Value = 16000
$ws.Range("O68").Value = 17000
$ws.Range("P68").Value = 16500
$ws.Range("Q68").Value = "$/caja 15 kilos"
$ws.Range("R68").Value = "Provincia de Limarí"
$ws.Range("S68").Value = 1100
$ws.Range("T68").Value = 15
